{"js": "// Update the \"Averages:\" summary row of the 12-month ACER table:\n//  - Label cell: \"Averages:\" -> \"Weighted Averages:\"\n//  - 12-month Survival average: \"39.07\" -> \"34.65\"\n//  - 12-Month % Dead average: \"45.26\" -> \"\" (cleared)\n//  - 12-Month % Missing average: \"15.67\" -> \"\" (cleared)\n//  - 12-Month % Cluster Survival average: \"58.58\" -> \"51.52\"\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in the document.\");\n}\n\n// The document is a single large table; find the summary row whose first\n// cell reads \"Averages:\" (trimmed) rather than assuming a fixed row index.\nlet targetRow = null;\n\nfor (const table of tables.items) {\n  const rows = table.rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  for (const row of rows.items) {\n    const cells = row.cells;\n    cells.load(\"items\");\n    await context.sync();\n\n    if (cells.items.length === 0) {\n      continue;\n    }\n\n    const firstCell = cells.items[0];\n    firstCell.load(\"value\");\n    await context.sync();\n\n    if ((firstCell.value || \"\").trim() === \"Averages:\") {\n      targetRow = row;\n      break;\n    }\n  }\n\n  if (targetRow) {\n    break;\n  }\n}\n\nif (!targetRow) {\n  throw new Error('Could not find the \"Averages:\" row.');\n}\n\nconst cells = targetRow.cells;\ncells.load(\"items\");\nawait context.sync();\n\n// Column order: Site ID(label) | Date | # Outplanted | Survival | % Dead | % Missing | % Cluster Survival\ncells.items[0].value = \"Weighted Averages:\";\ncells.items[3].value = \"34.65\";\ncells.items[4].value = \"\";\ncells.items[5].value = \"\";\ncells.items[6].value = \"51.52\";\n\nawait context.sync();\n", "ps1": "# Update the \"Averages:\" summary row of the 12-month ACER table:\n#  - Label cell: \"Averages:\" -> \"Weighted Averages:\"\n#  - 12-month Survival average: \"39.07\" -> \"34.65\"\n#  - 12-Month % Dead average: \"45.26\" -> \"\" (cleared)\n#  - 12-Month % Missing average: \"15.67\" -> \"\" (cleared)\n#  - 12-Month % Cluster Survival average: \"58.58\" -> \"51.52\"\n\n$d = $word.ActiveDocument\n\nif ($d.Tables.Count -eq 0) {\n    throw \"No table found in the document.\"\n}\n\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n\n$targetRow = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    $labelText = $table.Cell($r, 1).Range.Text\n    $labelText = ($labelText -replace \"[\\r\\a]\", \"\").Trim()\n    if ($labelText -eq \"Averages:\") {\n        $targetRow = $r\n        break\n    }\n}\n\nif ($targetRow -eq 0) {\n    throw 'Could not find the \"Averages:\" row.'\n}\n\n# Column order: Site ID(label) | Date | # Outplanted | Survival | % Dead | % Missing | % Cluster Survival\n$table.Cell($targetRow, 1).Range.Text = \"Weighted Averages:\"\n$table.Cell($targetRow, 4).Range.Text = \"34.65\"\n$table.Cell($targetRow, 5).Range.Text = \"\"\n$table.Cell($targetRow, 6).Range.Text = \"\"\n$table.Cell($targetRow, 7).Range.Text = \"51.52\"\n"}
